# OpenTBS 1.9.1 beta - debug parameter unique
#
# 1) The chart on slide 3 gets new (unique) category/value axis IDs when
#    OpenTBS merges/regenerates the chart XML. The PowerPoint object model
#    does not expose chart axis IDs for editing (Axis has no AxisId/Id
#    automation member), so that internal templating-library detail can't
#    be reproduced through COM - it is skipped here.
# 2) Refresh the cached "date updates automatically" placeholder text
#    (master + all 11 layouts) from 25/04/2013 to 22/07/2014.
# 3) Tidy the stray trailing paragraph-end run properties on the
#    "Merging a chart" slide title.

$p = $ppt.ActivePresentation

function Set-DatePlaceholderText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)

        $isDatePlaceholder = $false
        if ($sh.HasTextFrame) {
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
        }

        if ($isDatePlaceholder) {
            $tr = $sh.TextFrame.TextRange
            $tr.Delete()
            $tr.Text = $newText
        }
    }
}

$newDate = "22/07/2014"

# Slide master footer/date placeholder.
$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes $newDate

# Every slide layout has its own (inherited) date placeholder shape too.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes $newDate
}

# Slide 3 ("Merging a chart" title): drop the redundant endParaRPr that
# trails the single run in the title paragraph.
$slide3 = $p.Slides.Item(3)
$titleShape = $slide3.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$titleText = $titleRange.Text
$titleRange.Delete()
$titleRange.Text = $titleText
